{"js": "// Update the date line and the two-digit/one-digit division problems\n// in the table, per the commit's regenerated worksheet content.\n// Each (oldText -> newText) pair below corresponds 1:1 to a single\n// unique run in the original document, so every old value is looked\n// up and all matches collected BEFORE any text is written, which\n// keeps the pass safe even though some new values equal other old\n// values (e.g. \"30\u00f76=\" becomes \"65\u00f74=\", while the original \"65\u00f74=\"\n// cell becomes \"54\u00f74=\").\nconst replacements = [\n  [\"2024-12-18 Wednesday\", \"2024-12-19 Thursday\"],\n  [\"88\u00f79=\", \"17\u00f76=\"],\n  [\"15\u00f75=\", \"76\u00f72=\"],\n  [\"65\u00f74=\", \"54\u00f74=\"],\n  [\"79\u00f75=\", \"48\u00f72=\"],\n  [\"59\u00f72=\", \"95\u00f79=\"],\n  [\"41\u00f77=\", \"54\u00f73=\"],\n  [\"44\u00f74=\", \"34\u00f76=\"],\n  [\"13\u00f72=\", \"84\u00f78=\"],\n  [\"54\u00f78=\", \"44\u00f73=\"],\n  [\"94\u00f77=\", \"93\u00f78=\"],\n  [\"77\u00f73=\", \"23\u00f73=\"],\n  [\"39\u00f78=\", \"70\u00f75=\"],\n  [\"18\u00f73=\", \"40\u00f72=\"],\n  [\"41\u00f74=\", \"42\u00f77=\"],\n  [\"99\u00f78=\", \"23\u00f76=\"],\n  [\"90\u00f78=\", \"38\u00f79=\"],\n  [\"30\u00f76=\", \"65\u00f74=\"],\n  [\"19\u00f79=\", \"21\u00f73=\"],\n  [\"71\u00f77=\", \"94\u00f73=\"],\n  [\"61\u00f79=\", \"88\u00f72=\"],\n  [\"43\u00f75=\", \"41\u00f73=\"],\n  [\"66\u00f79=\", \"31\u00f76=\"],\n  [\"79\u00f73=\", \"56\u00f73=\"],\n  [\"46\u00f74=\", \"75\u00f75=\"],\n  [\"90\u00f74=\", \"27\u00f74=\"],\n];\n\nconst body = context.document.body;\n\n// First, resolve every search so we hold a Range for each occurrence\n// of the *original* text before any mutation happens.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: true })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Now perform the replacements using the ranges captured above.\nsearchResults.forEach((results, idx) => {\n  const [, newText] = replacements[idx];\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n});\n\nawait context.sync();\n", "ps1": "# Update the date line and the two-digit/one-digit division problems\n# in the table, per the commit's regenerated worksheet content.\n# Each (Old -> New) pair corresponds 1:1 to a single unique occurrence\n# of text in the document, so Find/Replace with MatchWholeWord is\n# unambiguous. Pairs are applied in document order, each keyed off the\n# *original* text, which keeps this correct even though some New\n# values equal other cells' Old values (e.g. \"30\u00f76=\" becomes \"65\u00f74=\",\n# while the original \"65\u00f74=\" cell becomes \"54\u00f74=\") \u2014 by the time we\n# search for \"65\u00f74=\" as a New target it no longer exists verbatim as\n# a find target for an earlier pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-12-18 Wednesday\"; New = \"2024-12-19 Thursday\" },\n    @{ Old = \"88\u00f79=\"; New = \"17\u00f76=\" },\n    @{ Old = \"15\u00f75=\"; New = \"76\u00f72=\" },\n    @{ Old = \"65\u00f74=\"; New = \"54\u00f74=\" },\n    @{ Old = \"79\u00f75=\"; New = \"48\u00f72=\" },\n    @{ Old = \"59\u00f72=\"; New = \"95\u00f79=\" },\n    @{ Old = \"41\u00f77=\"; New = \"54\u00f73=\" },\n    @{ Old = \"44\u00f74=\"; New = \"34\u00f76=\" },\n    @{ Old = \"13\u00f72=\"; New = \"84\u00f78=\" },\n    @{ Old = \"54\u00f78=\"; New = \"44\u00f73=\" },\n    @{ Old = \"94\u00f77=\"; New = \"93\u00f78=\" },\n    @{ Old = \"77\u00f73=\"; New = \"23\u00f73=\" },\n    @{ Old = \"39\u00f78=\"; New = \"70\u00f75=\" },\n    @{ Old = \"18\u00f73=\"; New = \"40\u00f72=\" },\n    @{ Old = \"41\u00f74=\"; New = \"42\u00f77=\" },\n    @{ Old = \"99\u00f78=\"; New = \"23\u00f76=\" },\n    @{ Old = \"90\u00f78=\"; New = \"38\u00f79=\" },\n    @{ Old = \"30\u00f76=\"; New = \"65\u00f74=\" },\n    @{ Old = \"19\u00f79=\"; New = \"21\u00f73=\" },\n    @{ Old = \"71\u00f77=\"; New = \"94\u00f73=\" },\n    @{ Old = \"61\u00f79=\"; New = \"88\u00f72=\" },\n    @{ Old = \"43\u00f75=\"; New = \"41\u00f73=\" },\n    @{ Old = \"66\u00f79=\"; New = \"31\u00f76=\" },\n    @{ Old = \"79\u00f73=\"; New = \"56\u00f73=\" },\n    @{ Old = \"46\u00f74=\"; New = \"75\u00f75=\" },\n    @{ Old = \"90\u00f74=\"; New = \"27\u00f74=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
